$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns (rows 2-10)
$values = @(
    @(5, 5),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("I$row").Value = $values[$i][0]
    $ws.Range("J$row").Value = $values[$i][1]
}
